$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the two empty "section header" rows (situação do domicílio / grandes regiões)
# which had a label in column A but no data. Delete the lower one first so the
# row index of the upper one does not shift before it is removed.
$ws.Rows(8).Delete()
$ws.Rows(5).Delete()

# Rename the "unnamed: 1_level_1" header in B2 to "total"
$ws.Range("B2").Value = "total"
